$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (existing invoice line) with new values ---
# GL_Date / Invoice_Date are stored as literal text (not real dates) in this
# workbook, so a leading apostrophe forces text entry instead of Excel's
# auto date-recognition.
$ws.Range("C2").Value = "'08/04/25"
$ws.Range("E2").Value = "JONSUP"
$ws.Range("G2").Value = "110-S101125942.001"
$ws.Range("I2").Value = "'08/04/25"
$ws.Range("J2").Value = 107.3
$ws.Range("R2").Value = "Ricky's Truck Stock"
$ws.Range("T2").Value = 5260
$ws.Range("U2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("W2").ClearContents()
$ws.Range("AA2").Value = "'2025"
$ws.Range("AB2").Value = "!Service Material"
$ws.Range("BF2").ClearContents()

# --- Remove row 3 entirely (its data has been merged into row 2 above) ---
$ws.Rows(3).Delete()
